$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.111.83"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.48%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.904.09"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.28%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "566.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.75"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.47%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.902.41"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.31%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.499"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.96"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.21%  "
$ws.Range("E11").Value = "  -1.28%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.431"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.27%  "
$ws.Range("E13").Value = "  +0.41%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.48"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.43%  "
$ws.Range("E15").Value = "  +0.41%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.386.98"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.049.64"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.43%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.916.49"
$ws.Range("D18").Style = "Normal"
$ws.Range("E19").Value = "  -0.69%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "429.69"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.01"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.41%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.651"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.12%  "
$ws.Range("E23").Value = "  -1.25%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.54"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.95"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.88%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.06"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.88%  "
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("E28").Value = "  -3.47%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000112"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.34%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.92"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.17%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.50"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.29%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.41%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.12%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.65"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.08%  "
$ws.Range("E35").Value = "  -3.69%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.956"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.00%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.37"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.50%  "
$ws.Range("E38").Value = "  -3.38%  "
$ws.Range("E39").Value = "  -0.55%  "
$ws.Range("E40").Value = "  -5.45%  "
$ws.Range("E41").Value = "  -1.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.12"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.88%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "40.64"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.14%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.713.66"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.54%  "
$ws.Range("E45").Value = "  -2.75%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "133.20"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.82%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0336"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.11%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "344.94"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.34%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000217"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +13.11%  "
$ws.Range("E51").Value = "  -0.78%  "
